$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = (Get-Date -Year 2018 -Month 2 -Day 6 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B11").Value = "absent du au décès de ma grand-maman"
$ws.Range("D11").Value = "1h30"

$ws.Range("A12").Value = (Get-Date -Year 2018 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B12").Value = "recherche d'une methode de cryptage du mot de passe"
$ws.Range("D12").Value = "1h30"

$ws.Range("A13").Value = (Get-Date -Year 2018 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B13").Value = "création de la classe cryptage"
$ws.Range("D13").Value = "2h15"

$ws.Range("D14").Select()
